$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows that were fully removed in the edit (delete higher index first)
$ws.Rows.Item(28).Delete()   # old "SC 92" row
$ws.Rows.Item(26).Delete()   # old "RM 232" row

# Apply the remaining per-cell value changes (post row-shift row numbers)
$ws.Range("D3").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("F6").Value = 16.43
$ws.Range("E8").Value = -6.6
$ws.Range("E10").Value = -6.1
$ws.Range("F11").Value = 17.65
$ws.Range("E12").ClearContents()
$ws.Range("F12").ClearContents()
$ws.Range("F13").Value = 17.1
$ws.Range("E15").Value = -8.4
$ws.Range("F17").ClearContents()
$ws.Range("E18").ClearContents()
$ws.Range("E19").ClearContents()
$ws.Range("F19").ClearContents()
$ws.Range("E25").Value = -7.1
$ws.Range("F25").Value = 16.6
$ws.Range("B26").Value = -20.2
$ws.Range("B27").ClearContents()
$ws.Range("E29").ClearContents()
$ws.Range("F31").ClearContents()
$ws.Range("F32").ClearContents()
$ws.Range("B33").Value = -19.5
$ws.Range("D33").Value = -14.1
